$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header names ---
# ExtractionType and SamplePortion swap positions (H1/I1)
$ws.Range("H1").Value = "ExtractionType"
$ws.Range("I1").Value = "SamplePortion"

# --- Row 2: type annotations ---
$ws.Range("A2").Value = "#string"
$ws.Range("B2").Value = "#string"
$ws.Range("C2").Value = "#date"
$ws.Range("D2").Value = "#string"
$ws.Range("E2").Value = "#string"
$ws.Range("F2").Value = "#string"
$ws.Range("G2").Value = "#string"
$ws.Range("H2").Value = "#string"
$ws.Range("I2").Value = "#float,  unit:µlormg"

# --- Row 3: French description/enum labels (new row) ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#TypeExtraction"
$ws.Range("I3").Value = "#PriseEssai"
